$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1025

# Row 3
$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 11 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21625
$ws.Range("S3").Value = 1081

# Row 4
$ws.Range("D4").Value = (Get-Date -Year 2021 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("S4").Value = 825

# Row 5
$ws.Range("D5").Value = (Get-Date -Year 2022 -Month 8 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 1075

# Row 6
$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1075

# Row 7
$ws.Range("D7").Value = (Get-Date -Year 2022 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22250
$ws.Range("S7").Value = 1112

# Row 8
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 725

# Row 9
$ws.Range("D9").Value = (Get-Date -Year 2023 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 24333
$ws.Range("S9").Value = 1217

# Row 10
$ws.Range("D10").Value = (Get-Date -Year 2022 -Month 5 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("S10").Value = 1025

# Row 11
$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("S11").Value = 825
